$wb = $excel.ActiveWorkbook

# "both_name_blank_after_mi_name" already demonstrates a multiindex whose
# second level shows just one blank (repeated-label) cell. Duplicate it to
# create a sibling fixture, "both_name_multiple_blanks", where every repeated
# row in that level is blanked out instead of only the first one, and place
# it right before "both_name_skiprows".
$src = $wb.Worksheets.Item("both_name_blank_after_mi_name")
$dst = $wb.Worksheets.Item("both_name_skiprows")

$src.Copy($dst) | Out-Null

# Restore the source sheet's selection over its full data range.
$src.Activate() | Out-Null
$src.Range("A1:F7").Select() | Out-Null

$new = $wb.Worksheets.Item("both_name_blank_after_mi_name (2)")
$new.Name = "both_name_multiple_blanks"

# Clear the remaining repeated index labels in column B so every row of the
# second index level is blank (multiple blanks), not just the first one.
$new.Range("B5").ClearContents() | Out-Null
$new.Range("B6").ClearContents() | Out-Null
$new.Range("B7").ClearContents() | Out-Null

# Make the new sheet the active tab, with the cursor left outside the data
# (mirroring the workbook's previous "parked" selection state).
$new.Activate() | Out-Null
